# Update res_bus vm_pu results for Case_5_122 (380 kV case)
# Columns: B,C,D,E,F,I,J,K,L,M for rows 2-25 (bus voltage magnitudes per unit)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.028784482558156
$ws.Cells.Item(2, 4).Value = 1.031106560225713
$ws.Cells.Item(2, 5).Value = 1.037373818706158
$ws.Cells.Item(2, 6).Value = 1.045197633122813
$ws.Cells.Item(2, 9).Value = 1.030645842577517
$ws.Cells.Item(2, 10).Value = 1.033934879726482
$ws.Cells.Item(2, 11).Value = 1.03391563697029
$ws.Cells.Item(2, 12).Value = 1.040164886640484
$ws.Cells.Item(2, 13).Value = 1.047966548433424
$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.029788170279179
$ws.Cells.Item(3, 4).Value = 1.031822661745746
$ws.Cells.Item(3, 5).Value = 1.038300522135754
$ws.Cells.Item(3, 6).Value = 1.046291838035138
$ws.Cells.Item(3, 9).Value = 1.030808273457226
$ws.Cells.Item(3, 10).Value = 1.034579045482416
$ws.Cells.Item(3, 11).Value = 1.034440480933844
$ws.Cells.Item(3, 12).Value = 1.040901074917713
$ws.Cells.Item(3, 13).Value = 1.048871408493036
$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.030437600828379
$ws.Cells.Item(4, 4).Value = 1.032285477444552
$ws.Cells.Item(4, 5).Value = 1.038900590218198
$ws.Cells.Item(4, 6).Value = 1.047000598492137
$ws.Cells.Item(4, 9).Value = 1.030911346223459
$ws.Cells.Item(4, 10).Value = 1.034995266394261
$ws.Cells.Item(4, 11).Value = 1.034778872226354
$ws.Cells.Item(4, 12).Value = 1.041377241684063
$ws.Cells.Item(4, 13).Value = 1.049457058461462
$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.030710615713436
$ws.Cells.Item(5, 4).Value = 1.032479912409638
$ws.Cells.Item(5, 5).Value = 1.039152960796598
$ws.Cells.Item(5, 6).Value = 1.047298737200053
$ws.Cells.Item(5, 9).Value = 1.030954191198927
$ws.Cells.Item(5, 10).Value = 1.035170102035964
$ws.Cells.Item(5, 11).Value = 1.034920839546412
$ws.Cells.Item(5, 12).Value = 1.041577374598863
$ws.Cells.Item(5, 13).Value = 1.049703299845839
$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.030756455769146
$ws.Cells.Item(6, 4).Value = 1.032512551078431
$ws.Cells.Item(6, 5).Value = 1.039195340884218
$ws.Cells.Item(6, 6).Value = 1.047348806294415
$ws.Cells.Item(6, 9).Value = 1.030961356493581
$ws.Cells.Item(6, 10).Value = 1.0351994492692
$ws.Cells.Item(6, 11).Value = 1.03494465931475
$ws.Cells.Item(6, 12).Value = 1.041610974987312
$ws.Cells.Item(6, 13).Value = 1.049744646842831
$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.030441248889518
$ws.Cells.Item(7, 4).Value = 1.032288076017131
$ws.Cells.Item(7, 5).Value = 1.038903962006546
$ws.Cells.Item(7, 6).Value = 1.04700458154603
$ws.Cells.Item(7, 9).Value = 1.030911920633764
$ws.Cells.Item(7, 10).Value = 1.034997603121807
$ws.Cells.Item(7, 11).Value = 1.034780770350541
$ws.Cells.Item(7, 12).Value = 1.041379916057647
$ws.Cells.Item(7, 13).Value = 1.049460348617486
$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.029123687990597
$ws.Cells.Item(8, 4).Value = 1.031348683273746
$ws.Cells.Item(8, 5).Value = 1.037686913194885
$ws.Cells.Item(8, 6).Value = 1.04556727191969
$ws.Cells.Item(8, 9).Value = 1.030701156989108
$ws.Cells.Item(8, 10).Value = 1.03415270189175
$ws.Cells.Item(8, 11).Value = 1.03409326211532
$ws.Cells.Item(8, 12).Value = 1.04041372519336
$ws.Cells.Item(8, 13).Value = 1.048272319757126
$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.026801818607666
$ws.Cells.Item(9, 4).Value = 1.029689190448959
$ws.Cells.Item(9, 5).Value = 1.035545635706725
$ws.Cells.Item(9, 6).Value = 1.043040215130842
$ws.Cells.Item(9, 9).Value = 1.030314237353928
$ws.Cells.Item(9, 10).Value = 1.032659328914762
$ws.Cells.Item(9, 11).Value = 1.032872493271132
$ws.Cells.Item(9, 12).Value = 1.038709696391907
$ws.Cells.Item(9, 13).Value = 1.046179996570819
$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.025253814436939
$ws.Cells.Item(10, 4).Value = 1.028580126586299
$ws.Cells.Item(10, 5).Value = 1.034120389714617
$ws.Cells.Item(10, 6).Value = 1.04135935013802
$ws.Cells.Item(10, 9).Value = 1.030045886450925
$ws.Cells.Item(10, 10).Value = 1.031660721950384
$ws.Cells.Item(10, 11).Value = 1.032052446663378
$ws.Cells.Item(10, 12).Value = 1.037572717384434
$ws.Cells.Item(10, 13).Value = 1.044785909037275
$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.024583492176191
$ws.Cells.Item(11, 4).Value = 1.028099253690401
$ws.Cells.Item(11, 5).Value = 1.033503789475113
$ws.Cells.Item(11, 6).Value = 1.040632433129845
$ws.Cells.Item(11, 9).Value = 1.0299272267649
$ws.Cells.Item(11, 10).Value = 1.031227603061432
$ws.Cells.Item(11, 11).Value = 1.031695895887648
$ws.Cells.Item(11, 12).Value = 1.037080171633254
$ws.Cells.Item(11, 13).Value = 1.044182448012638
$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.024334500767586
$ws.Cells.Item(12, 4).Value = 1.027920540753367
$ws.Cells.Item(12, 5).Value = 1.033274838493953
$ws.Cells.Item(12, 6).Value = 1.040362560637946
$ws.Cells.Item(12, 9).Value = 1.029882781895309
$ws.Cells.Item(12, 10).Value = 1.031066616368233
$ws.Cells.Item(12, 11).Value = 1.031563237590667
$ws.Cells.Item(12, 12).Value = 1.036897184641552
$ws.Cells.Item(12, 13).Value = 1.043958324354161
$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.024387910444246
$ws.Cells.Item(13, 4).Value = 1.027958879583721
$ws.Cells.Item(13, 5).Value = 1.033323945545775
$ws.Cells.Item(13, 6).Value = 1.040420443005336
$ws.Cells.Item(13, 9).Value = 1.029892332177841
$ws.Cells.Item(13, 10).Value = 1.031101153397776
$ws.Cells.Item(13, 11).Value = 1.031591703166298
$ws.Cells.Item(13, 12).Value = 1.036936437484557
$ws.Cells.Item(13, 13).Value = 1.044006398337861
$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.024562910547347
$ws.Cells.Item(14, 4).Value = 1.02808448316252
$ws.Cells.Item(14, 5).Value = 1.033484862643568
$ws.Cells.Item(14, 6).Value = 1.040610122591801
$ws.Cells.Item(14, 9).Value = 1.029923560470435
$ws.Cells.Item(14, 10).Value = 1.031214298032477
$ws.Cells.Item(14, 11).Value = 1.031684934784614
$ws.Cells.Item(14, 12).Value = 1.03706504655373
$ws.Cells.Item(14, 13).Value = 1.044163921294784
$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.024670733375118
$ws.Cells.Item(15, 4).Value = 1.028161859068883
$ws.Cells.Item(15, 5).Value = 1.033584019835006
$ws.Cells.Item(15, 6).Value = 1.04072700858945
$ws.Cells.Item(15, 9).Value = 1.029942752320873
$ws.Cells.Item(15, 10).Value = 1.031283996003652
$ws.Cells.Item(15, 11).Value = 1.031742348800371
$ws.Cells.Item(15, 12).Value = 1.037144282410615
$ws.Cells.Item(15, 13).Value = 1.044260980193173
$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.025298300970698
$ws.Cells.Item(16, 4).Value = 1.028612027115221
$ws.Cells.Item(16, 5).Value = 1.034161322882121
$ws.Cells.Item(16, 6).Value = 1.041407612398603
$ws.Cells.Item(16, 9).Value = 1.030053709675769
$ws.Cells.Item(16, 10).Value = 1.031689451587863
$ws.Cells.Item(16, 11).Value = 1.032076078951205
$ws.Cells.Item(16, 12).Value = 1.037605401292371
$ws.Cells.Item(16, 13).Value = 1.044825962752996
$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.025691950566819
$ws.Cells.Item(17, 4).Value = 1.028894234723932
$ws.Cells.Item(17, 5).Value = 1.034523594991379
$ws.Cells.Item(17, 6).Value = 1.041834780692858
$ws.Cells.Item(17, 9).Value = 1.030122651511215
$ws.Cells.Item(17, 10).Value = 1.031943591920998
$ws.Cells.Item(17, 11).Value = 1.03228502725792
$ws.Cells.Item(17, 12).Value = 1.037894588696133
$ws.Cells.Item(17, 13).Value = 1.045180412008272
$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.025921557180135
$ws.Cells.Item(18, 4).Value = 1.02905877974382
$ws.Cells.Item(18, 5).Value = 1.034734954313334
$ws.Cells.Item(18, 6).Value = 1.042084028378346
$ws.Cells.Item(18, 9).Value = 1.030162626539124
$ws.Cells.Item(18, 10).Value = 1.032091758614587
$ws.Cells.Item(18, 11).Value = 1.03240676185712
$ws.Cells.Item(18, 12).Value = 1.038063244948771
$ws.Cells.Item(18, 13).Value = 1.045387174383568
$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.025999846669453
$ws.Cells.Item(19, 4).Value = 1.029114874814111
$ws.Cells.Item(19, 5).Value = 1.034807031204786
$ws.Cells.Item(19, 6).Value = 1.042169030253943
$ws.Cells.Item(19, 9).Value = 1.030176216677733
$ws.Cells.Item(19, 10).Value = 1.032142267900379
$ws.Cells.Item(19, 11).Value = 1.032448246203135
$ws.Cells.Item(19, 12).Value = 1.038120748650061
$ws.Cells.Item(19, 13).Value = 1.04545767807247
$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.025649715952317
$ws.Cells.Item(20, 4).Value = 1.028863962914274
$ws.Cells.Item(20, 5).Value = 1.034484721226317
$ws.Cells.Item(20, 6).Value = 1.041788940509296
$ws.Cells.Item(20, 9).Value = 1.030115279274865
$ws.Cells.Item(20, 10).Value = 1.031916332208654
$ws.Cells.Item(20, 11).Value = 1.032262623704202
$ws.Cells.Item(20, 12).Value = 1.037863563883154
$ws.Cells.Item(20, 13).Value = 1.045142381080423
$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.024511377485183
$ws.Cells.Item(21, 4).Value = 1.028047498663144
$ws.Cells.Item(21, 5).Value = 1.033437474306325
$ws.Cells.Item(21, 6).Value = 1.040554262900752
$ws.Cells.Item(21, 9).Value = 1.029914374709253
$ws.Cells.Item(21, 10).Value = 1.031180982709297
$ws.Cells.Item(21, 11).Value = 1.031657486449731
$ws.Cells.Item(21, 12).Value = 1.037027175311549
$ws.Cells.Item(21, 13).Value = 1.04411753391663
$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.023795636523233
$ws.Cells.Item(22, 4).Value = 1.027533604085711
$ws.Cells.Item(22, 5).Value = 1.032779502162881
$ws.Cells.Item(22, 6).Value = 1.039778764027915
$ws.Cells.Item(22, 9).Value = 1.029785921168039
$ws.Cells.Item(22, 10).Value = 1.030718020130533
$ws.Cells.Item(22, 11).Value = 1.031275743215042
$ws.Cells.Item(22, 12).Value = 1.036501110612731
$ws.Cells.Item(22, 13).Value = 1.043473337688232
$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.024175066449704
$ws.Cells.Item(23, 4).Value = 1.027806081244407
$ws.Cells.Item(23, 5).Value = 1.033128260567862
$ws.Cells.Item(23, 6).Value = 1.040189795473325
$ws.Cells.Item(23, 9).Value = 1.029854219204063
$ws.Cells.Item(23, 10).Value = 1.030963503896917
$ws.Cells.Item(23, 11).Value = 1.031478232633884
$ws.Cells.Item(23, 12).Value = 1.036780005708352
$ws.Cells.Item(23, 13).Value = 1.043814822450355
$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.025668799964337
$ws.Cells.Item(24, 4).Value = 1.028877641633324
$ws.Cells.Item(24, 5).Value = 1.03450228644634
$ws.Cells.Item(24, 6).Value = 1.041809653444177
$ws.Cells.Item(24, 9).Value = 1.030118611205219
$ws.Cells.Item(24, 10).Value = 1.03192864991263
$ws.Cells.Item(24, 11).Value = 1.03227274734271
$ws.Cells.Item(24, 12).Value = 1.037877582728211
$ws.Cells.Item(24, 13).Value = 1.04515956556363
$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.02740209470015
$ws.Cells.Item(25, 4).Value = 1.030118695850171
$ws.Cells.Item(25, 5).Value = 1.036098809767904
$ws.Cells.Item(25, 6).Value = 1.043692845305579
$ws.Cells.Item(25, 9).Value = 1.030416101758365
$ws.Cells.Item(25, 10).Value = 1.033045937233087
$ws.Cells.Item(25, 11).Value = 1.033189187557989
$ws.Cells.Item(25, 12).Value = 1.039150400273363
$ws.Cells.Item(25, 13).Value = 1.049460348617486

Write-Host "vm_pu values updated for 380 kV case"
